$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 45
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 3).Value = 45177
}
